$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I39").Value = 209.6
$ws.Range("K39").Value = 628.8
$ws.Range("M39").Value = -332.8
$ws.Range("H39").Value = 264
$ws.Range("N112").Value = -10427.72
$ws.Range("H112").Value = 2737.24
$ws.Range("L112").Value = 8211.719999999999
$ws.Range("J112").Value = 2737.24
$ws.Range("N136").Value = -55200
$ws.Range("L136").Value = 45000
$ws.Range("H136").Value = 46250
$ws.Range("J136").Value = 45000
$ws.Range("J138").Value = 3642.3462
$ws.Range("H138").Value = 3623.2942
$ws.Range("N138").Value = -21207.0386
$ws.Range("L138").Value = 10927.0386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3639192.8
$ws.Range("K32").Value = 4084308
$ws.Range("M32").Value = -4084021
$ws.Range("I32").Value = 4084308
$ws.Range("K61").Value = 5062
$ws.Range("I61").Value = 5062
$ws.Range("H61").Value = 7437
$ws.Range("M61").Value = -4850
$ws.Range("K110").Value = 39231332
$ws.Range("H110").Value = 20402448
$ws.Range("M110").Value = -39229287
$ws.Range("I110").Value = 39231332
$ws.Range("H132").Value = 8768.571
$ws.Range("I132").Value = 5499.6924
$ws.Range("K132").Value = 16499.0772
$ws.Range("M132").Value = -13969.0772
$ws.Range("K136").Value = 15186
$ws.Range("I136").Value = 5062
$ws.Range("M136").Value = -12636
$ws.Range("H136").Value = 7437

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2026.2142
$ws.Range("N22").Value = -2563.3333
$ws.Range("J22").Value = 2217.3333
$ws.Range("L22").Value = 2217.3333
$ws.Range("H94").Value = 597.2105
$ws.Range("M94").Value = -110.4815
$ws.Range("J94").Value = 684.9091
$ws.Range("K94").Value = 561.4815
$ws.Range("L94").Value = 684.9091
$ws.Range("I94").Value = 561.4815
$ws.Range("N94").Value = -1586.9091
$ws.Range("H132").Value = 102499.5
$ws.Range("L132").Value = 102499.5
$ws.Range("J132").Value = 102499.5
$ws.Range("N132").Value = -112619.5
$ws.Range("N134").Value = -31010.7
$ws.Range("J134").Value = 8646.9
$ws.Range("L134").Value = 25940.7
$ws.Range("M134").Value = -6236.400000000001
$ws.Range("I134").Value = 2923.8
$ws.Range("K134").Value = 8771.400000000001
$ws.Range("H134").Value = 6194.143

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 422.7143
$ws.Range("I22").Value = 465.27274
$ws.Range("M22").Value = -115.27274
$ws.Range("K22").Value = 465.27274
$ws.Range("H31").Value = 4237.0264
$ws.Range("K31").Value = 2501.353
$ws.Range("J31").Value = 5642.095
$ws.Range("L31").Value = 5642.095
$ws.Range("N31").Value = -6232.095
$ws.Range("I31").Value = 2501.353
$ws.Range("M31").Value = -2206.353
$ws.Range("L34").Value = 5642.095
$ws.Range("N34").Value = -6046.095
$ws.Range("I34").Value = 2501.353
$ws.Range("H34").Value = 4237.0264
$ws.Range("J34").Value = 5642.095
$ws.Range("M34").Value = -2299.353
$ws.Range("K34").Value = 2501.353
$ws.Range("N122").Value = -13485.25
$ws.Range("J122").Value = 2861.75
$ws.Range("I122").Value = 170485.5
$ws.Range("K122").Value = 511456.5
$ws.Range("M122").Value = -509006.5
$ws.Range("H122").Value = 103436
$ws.Range("L122").Value = 8585.25
$ws.Range("H132").Value = 60285.137
$ws.Range("I132").Value = 5589.625
$ws.Range("K132").Value = 16768.875
$ws.Range("M132").Value = -14238.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 7165.2
$ws.Range("M2").Value = -1057
$ws.Range("N2").Value = -43217.2
$ws.Range("H2").Value = 4295.1177
$ws.Range("L2").Value = 42991.2
$ws.Range("I2").Value = 195
$ws.Range("K2").Value = 1170
$ws.Range("L37").Value = 219994.8
$ws.Range("H37").Value = 73331.60000000001
$ws.Range("J37").Value = 73331.60000000001
$ws.Range("N37").Value = -220218.8
$ws.Range("N51").Value = -8363.999899999999
$ws.Range("K51").Value = 4987.5
$ws.Range("H51").Value = 2013.4286
$ws.Range("L51").Value = 7443.999899999999
$ws.Range("I51").Value = 1662.5
$ws.Range("M51").Value = -4527.5
$ws.Range("J51").Value = 2481.3333
$ws.Range("I108").Value = 398.9
$ws.Range("H108").Value = 398.9
$ws.Range("M108").Value = 1683.3
$ws.Range("K108").Value = 1196.7
$ws.Range("H109").Value = 806
$ws.Range("M109").Value = -2213.2
$ws.Range("N109").Value = -2410
$ws.Range("J109").Value = 110
$ws.Range("K109").Value = 3253.2
$ws.Range("I109").Value = 1084.4
$ws.Range("L109").Value = 330
$ws.Range("K119").Value = 9361.7145
$ws.Range("H119").Value = 5794.5454
$ws.Range("I119").Value = 3120.5715
$ws.Range("M119").Value = -4523.7145
$ws.Range("M134").Value = -107139798
$ws.Range("I134").Value = 35714956
$ws.Range("K134").Value = 107144868
$ws.Range("H134").Value = 33334428

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L110").Value = 0
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("I114").Value = 0
$ws.Range("H114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("K114").Value = 0
$ws.Range("K126").Value = 6480.333
$ws.Range("H126").Value = 3246.25
$ws.Range("I126").Value = 2160.111
$ws.Range("L126").Value = 19513.9995
$ws.Range("M126").Value = -4010.333
$ws.Range("J126").Value = 6504.6665
$ws.Range("N126").Value = -24453.9995
$ws.Range("H132").Value = 9251.377
$ws.Range("I132").Value = 7074.3784
$ws.Range("K132").Value = 21223.1352
$ws.Range("M132").Value = -18693.1352

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2664.7551
$ws.Range("N22").Value = -3903.6924
$ws.Range("J22").Value = 3313.6924
$ws.Range("L22").Value = 3313.6924
$ws.Range("J27").Value = 3313.6924
$ws.Range("N27").Value = -3527.6924
$ws.Range("L27").Value = 3313.6924
$ws.Range("H27").Value = 2664.7551
$ws.Range("M33").Value = -39725
$ws.Range("H33").Value = 40015
$ws.Range("K33").Value = 40015
$ws.Range("I33").Value = 40015
$ws.Range("H40").Value = 35719560
$ws.Range("M40").Value = -41671980
$ws.Range("I40").Value = 41672116
$ws.Range("K40").Value = 41672116
$ws.Range("J55").Value = 1008.625
$ws.Range("K55").Value = 331.8
$ws.Range("N55").Value = -1354.625
$ws.Range("L55").Value = 1008.625
$ws.Range("I55").Value = 331.8
$ws.Range("M55").Value = -158.8
$ws.Range("H55").Value = 632.6111
$ws.Range("N122").Value = -23344.375
$ws.Range("J122").Value = 6148.125
$ws.Range("I122").Value = 29414574
$ws.Range("K122").Value = 88243722
$ws.Range("M122").Value = -88241272
$ws.Range("H122").Value = 20003878
$ws.Range("L122").Value = 18444.375
$ws.Range("H132").Value = 6479.533
$ws.Range("L132").Value = 30749.25
$ws.Range("J132").Value = 10249.75
$ws.Range("I132").Value = 5108.5454
$ws.Range("N132").Value = -35809.25
$ws.Range("K132").Value = 15325.6362
$ws.Range("M132").Value = -12795.6362
$ws.Range("K136").Value = 8367.999899999999
$ws.Range("I136").Value = 2789.3333
$ws.Range("N136").Value = -35474.25
$ws.Range("M136").Value = -5817.999899999999
$ws.Range("L136").Value = 30374.25
$ws.Range("H136").Value = 6241.294
$ws.Range("J136").Value = 10124.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K126").Value = 60608286
$ws.Range("H126").Value = 14431965
$ws.Range("I126").Value = 20202762
$ws.Range("L126").Value = 14919
$ws.Range("M126").Value = -60605816
$ws.Range("J126").Value = 4973
$ws.Range("N126").Value = -19859

Write-Output "applied 194 cell changes across 8 sheets"